{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(Word.RangeLocation.end);\n\nconst newParagraphsOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p/><w:p><w:r><w:t>My Milestone</w:t></w:r><w:r><w:t xml:space=\"preserve\"> 1</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">I put </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">overview and outline </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>is pushed on</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>github</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> also</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>my start code is on</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> G</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>tHub</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> repository.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">Organized </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">Setting up to </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">start </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">and </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>UI software framework</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> for React Native</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> and </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>learning  new</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> language</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> for the platform</w:t></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(newParagraphsOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$endPos = $d.Content.End\n$insertionRange = $d.Range($endPos, $endPos)\n\n$newParagraphsOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p/><w:p><w:r><w:t>My Milestone</w:t></w:r><w:r><w:t xml:space=\"preserve\"> 1</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr><w:r><w:t xml:space=\"preserve\">I put </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">overview and outline </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>is pushed on</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>github</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> also</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>my start code is on</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> G</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>tHub</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> repository.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">Organized </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">Setting up to </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">start </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">and </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\">learning </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>UI software framework</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> for React Native</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> and </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t>learning  new</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> language</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Lato\" w:hAnsi=\"Lato\"/><w:color w:val=\"000000\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/></w:rPr><w:t xml:space=\"preserve\"> for the platform</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$insertionRange.InsertXML($newParagraphsOoxml)\n"}
